$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Business Unit values in column H (rows 2 and 3) to the corrected,
# properly-capitalized business unit name.
$ws.Range("H2").Value = "Graha-Segara-Belawan"
$ws.Range("H3").Value = "Graha-Segara-Belawan"

# Remove the now-unneeded "Status"/contract column (column M), shifting all
# subsequent columns (Agama, No Kk, No Ktp, ... Gol Darah) one position to
# the left.
$ws.Columns("M").Delete()

# Reset the view/selection back to the top-left of the data.
$null = $ws.Range("H9").Select()
